$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8 (pushes existing rows 8-19 down to 9-20)
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new weekly data point
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 44533
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = 100112040
$ws.Cells.Item(8, 7).Value = "Cilantro"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 2000
$ws.Cells.Item(8, 12).Value = 2200
$ws.Cells.Item(8, 13).Value = 2100
$ws.Cells.Item(8, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(8, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(8, 16).Value = 2100
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = "Hortaliza"
